$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.104.28'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +0.24%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.837.99'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +0.55%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  +0.29%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.49'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.90%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6270'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -1.39%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.003'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +0.32%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07588'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +3.42%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2932'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +0.14%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '22.59'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -1.12%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07746'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +1.16%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.831.91'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +0.28%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.965'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -0.31%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6650'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +0.07%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.00001011'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +17.77%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '82.86'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +0.75%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.061'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -0.02%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '29.073.18'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +0.54%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '227.18'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +1.50%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.38'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.09%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.004'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +0.45%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.205'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +1.69%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.003'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +0.28%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '159.37'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.91%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.522'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +0.99%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1384'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +0.40%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.95'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +0.42%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.499'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -0.47%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.108'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +0.58%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.023'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +0.27%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.196'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -0.24%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.05256'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -0.75%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.844'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +0.90%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7343'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -0.90%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.139'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -1.06%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.699'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +2.17%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.245.04'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -3.40%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.766'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +0.93%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01787'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +0.31%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.356'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +0.48%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8980'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +0.37%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.002'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +0.39%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '102.26'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -0.33%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.973.84'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -0.25%  '
$ws.Range("E45").Value = '  -1.71%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '64.37'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +0.42%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5122'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -0.19%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4050'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +1.81%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.857'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +1.67%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05759'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -1.25%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.686'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +0.18%  '
